$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the country names in A24 / A25 (Brasil <-> Japon)
$ws.Range("A24").Value = "Japon"
$ws.Range("A25").Value = "Brasil"

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 19:46"

# Updated case numbers for several countries
$ws.Range("B7").Value = 22813
$ws.Range("C7").Value = 3430
$ws.Range("E7").Value = 22354

$ws.Range("B8").Value = 22084
$ws.Range("C8").Value = 2236
$ws.Range("E8").Value = 21792
$ws.Range("G8").Value = 15
$ws.Range("H8").Value = 83

$ws.Range("B10").Value = 14459
$ws.Range("C10").Value = 1847
$ws.Range("E10").Value = 12310
$ws.Range("F10").Value = 1525
$ws.Range("G10").Value = 112
$ws.Range("H10").Value = 562

$ws.Range("B12").Value = 6489
$ws.Range("C12").Value = 874
$ws.Range("E12").Value = 6400

$ws.Range("B13").Value = 5018
$ws.Range("C13").Value = 1035
$ws.Range("E13").Value = 4720

$ws.Range("B17").Value = 2128
$ws.Range("C17").Value = 169
$ws.Range("E17").Value = 2120

$ws.Range("B22").Value = 1145
$ws.Range("C22").Value = 58
$ws.Range("E22").Value = 1118

$ws.Range("B24").Value = 1046
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 215
$ws.Range("E24").Value = 795
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 36

$ws.Range("B25").Value = 1021
$ws.Range("C25").Value = 51
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 1001
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 18
